# OneDM conference schedule update (from 5/6 Session 2):
#  - Tue 13:00-15:00 slot (row 6): the "Bus/Tech 1" column (E6) becomes a
#    cancelled slot, highlighted in red like the other cancelled/holiday cells.
#  - Tue 13:00-15:00 slot (row 6): the "Bus 2" column (F6) gets an expanded
#    description.
#  - Thu 16:00-18:00 slot (row 8): the "Tech - SDF versions..." column (E8)
#    is replaced with a note continuing the Wednesday discussion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E6: "Bus/Tech 1" -> "Cancel", styled red (font + fill) like the other
# cancelled/holiday cells (column B).
$ws.Range("E6").Value = "Cancel"
$ws.Range("E6").Font.Color = 255
$ws.Range("E6").Interior.Color = 255

# E8: "Tech - SDF versions, language features" -> continuation note.
$ws.Range("E8").Value = "Continue Discussion from Wed. + SDF Specification including I/D"

# F6: "Bus 2 (David Mc.) Public statements, etc." -> expanded description.
$ws.Range("F6").Value = "Business/nontec (David Mc.) Public statements, etc. FAQ and explainers"

# Scroll/selection moved when the file was last saved.
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("F7").Select()
